# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" stats table: revised per-country counters
# (some countries overtake their neighbours in the ranking, swapping rows)
# and bump the "last updated" timestamp string in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 31 de Julio de 2020 a las 01:02"

# --- Row 4: Estados Unidos (in-place refresh) ---------------------------
$ws.Range("B4").Value = 4626656
$ws.Range("C4").Value = 58619
$ws.Range("D4").Value = 2276208
$ws.Range("E4").Value = 2195461
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1147
$ws.Range("H4").Value = 154987

# --- Row 5: Brasil (in-place refresh) ------------------------------------
$ws.Range("D5").Value = 1824095
$ws.Range("E5").Value = 694744

# --- Row 8: Sudafrica (in-place refresh) --------------------------------
$ws.Range("B8").Value = 482169
$ws.Range("C8").Value = 11046
$ws.Range("D8").Value = 309601
$ws.Range("E8").Value = 164756
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 315
$ws.Range("H8").Value = 7812

# --- Rows 15-16: Colombia overtakes Pakistan, swapping rank positions ---
$ws.Range("A15").Value = "Colombia"
$ws.Range("B15").Value = 286020
$ws.Range("C15").Value = 9965
$ws.Range("D15").Value = 148695
$ws.Range("E15").Value = 127515
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 356
$ws.Range("H15").Value = 9810

$ws.Range("A16").Value = "Pakistan"
$ws.Range("B16").Value = 277402
$ws.Range("C16").Value = 1114
$ws.Range("D16").Value = 246131
$ws.Range("E16").Value = 25347
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = 5924

# --- Row 23: Argentina (in-place refresh) -------------------------------
$ws.Range("B23").Value = 185373
$ws.Range("C23").Value = 6377
$ws.Range("D23").Value = 80596
$ws.Range("E23").Value = 101336
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 153
$ws.Range("H23").Value = 3441

# --- Row 50: Nigeria (in-place refresh) ---------------------------------
$ws.Range("B50").Value = 42689
$ws.Range("C50").Value = 481
$ws.Range("D50").Value = 19270
$ws.Range("E50").Value = 22541
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 5
$ws.Range("H50").Value = 878

# --- Row 58: Japon (in-place refresh) -----------------------------------
$ws.Range("B58").Value = 33049
$ws.Range("C58").Value = 1148
$ws.Range("D58").Value = 24179
$ws.Range("E58").Value = 7866
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 1004

# --- Rows 64-65: Uzbekistan overtakes Marruecos, swapping rank positions -
$ws.Range("A64").Value = "Uzbekistan"
$ws.Range("B64").Value = 23271
$ws.Range("C64").Value = 686
$ws.Range("D64").Value = 13680
$ws.Range("E64").Value = 9455
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 5
$ws.Range("H64").Value = 136

$ws.Range("A65").Value = "Marruecos"
$ws.Range("B65").Value = 23259
$ws.Range("C65").Value = 1046
$ws.Range("D65").Value = 17311
$ws.Range("E65").Value = 5602
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 12
$ws.Range("H65").Value = 346

# --- Row 86: Noruega (in-place refresh) ---------------------------------
$ws.Range("B86").Value = 9208
$ws.Range("C86").Value = 36
$ws.Range("D86").Value = 8752
$ws.Range("E86").Value = 201
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 255

# --- Row 139: Uruguay (in-place refresh) --------------------------------
$ws.Range("B139").Value = 1243
$ws.Range("C139").Value = 6
$ws.Range("D139").Value = 978
$ws.Range("E139").Value = 230
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 35

# --- Row 162: Bahamas (in-place refresh) --------------------------------
$ws.Range("B162").Value = 508
$ws.Range("C162").Value = 24
$ws.Range("D162").Value = 91
$ws.Range("E162").Value = 403
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 3
$ws.Range("H162").Value = 14
